$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G10").Value = 'esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G18").Value = 'yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G19").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G20").Value = 'marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G21").Value = 'yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G24").Value = 'majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G25").Value = 'backup@backdoor.com, gehanadel@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg'
$ws.Range("G32").Value = 'esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G40").Value = 'yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G41").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G42").Value = 'marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G46").Value = 'nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G54").Value = 'Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G58").Value = 'Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range("G62").Value = 'yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G63").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G64").Value = 'wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range("G65").Value = 'ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G68").Value = 'nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G76").Value = 'Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G80").Value = 'Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range("G81").Value = 'user@user.com, Walaa.h.ghanima@med.asu.edu.eg, enas.omran@med.asu.edu.eg'
$ws.Range("G84").Value = 'yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G85").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G86").Value = 'wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range("G90").Value = 'Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Range("G96").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G98").Value = 'Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G106").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G107").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G108").Value = 'marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G112").Value = 'Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Range("G118").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G120").Value = 'Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G128").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G129").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G130").Value = 'marinasorial@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G134").Value = 'majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G142").Value = 'esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
$ws.Range("G150").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G151").Value = 'marina_atef@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G156").Value = 'majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Range("G164").Value = 'esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
$ws.Range("G165").Value = 'Sarah.Mahdy@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg'
$ws.Range("G172").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G173").Value = 'marina_atef@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
